# Generate Report for Archive
#
# Two files (a47f2c1d-0b3f-4a3a-bda4-f8241ec2f013.md and
# cff7eb05-913d-4886-9496-e79eb5af5fa8.md) moved from "Ready for handoff"
# to "In Translation" status. Update the per-language status columns on the
# zh-cn and de-de sheets, as well as the summary status columns on the
# Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: rows 3 and 4 correspond to a47f2c1d-...md and
# cff7eb05-...md; columns E (zh-cn) and F (de-de) hold the status text.
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

# zh-cn / de-de detail sheets: rows 3 and 4 are the same two files;
# column C holds the Status value.
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
